$d = $word.ActiveDocument

$d.Content.Find.Execute(
    "php8.1-imap php8.1-redis php8.1-snmp",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "php8.1-mysql",
    2
)
